$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45172 to 45175 for rows 2-77
for ($r = 2; $r -le 77; $r++) {
    $ws.Cells.Item($r, 3).Value = 45175
}

# Add new row 78 with data
$ws.Range("A78").Value = "A 40897-2023"
$ws.Range("B78").Value = 45173
$ws.Range("C78").Value = 45175
$ws.Range("D78").Value = "ÖREBRO LÄN"
$ws.Range("E78").Value = "KUMLA"
$ws.Range("G78").Value = 4.7
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0

# Apply same styles as row 77 for B,C (date format) and R (wrap text)
$ws.Range("B78").NumberFormat = $ws.Range("B77").NumberFormat
$ws.Range("C78").NumberFormat = $ws.Range("C77").NumberFormat
$ws.Range("R78").WrapText = $true

# Row 77 picks up an explicit row height (matches the rest of the sheet)
$ws.Rows.Item(77).RowHeight = 15
